# PartsList_LidarBreakout.xlsx edits
# - Swapped USB connector (row 12) from TE Connectivity Micro-AB to Hirose Mini-B
# - Updated electrolytic cap (row 23) from 6.3V SMD to 25V radial part
# - Added 3 new debounce-related parts (MOSFET + two caps) as new rows 32-34
# - Moved the MAX/BOARD total down to row 41 and widened the SUM range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: USB connector swap (Micro-AB -> Mini-B) ---
$ws.Range("A12").Value2 = "H2960CT-ND"
$ws.Range("D12").Value2 = 0.99
$ws.Range("E12").Value2 = "UX60-MB-5S8"
$ws.Range("F12").Value2 = "Hirose Electric Co Ltd"
$ws.Range("G12").Value2 = "CONN RCPT USB2.0 MINI B SMD R/A"

# --- Row 23: electrolytic capacitor swap (6.3V SMD -> 25V radial) ---
$ws.Range("A23").Value2 = "P15367CT-ND"
$ws.Range("E23").Value2 = "EEU-FR1E681LB"
$ws.Range("G23").Value2 = "CAP ALUM 680UF 20% 25V RADIAL"

# --- Insert 3 new rows for debounce caps/mosfet before the old trailing blank rows ---
$ws.Rows("32:34").Insert()

# Row 32: MOSFET for debounce
$ws.Range("A32").Value2 = "BSS84PH6433XTMA1CT-ND"
$ws.Range("B32").Value2 = 3
$ws.Range("C32").Value2 = 1
$ws.Range("D32").Value2 = 0.39
$ws.Range("E32").Value2 = "BSS84PH6433XTMA1"
$ws.Range("F32").Value2 = "Infineon Technologies"
$ws.Range("G32").Value2 = "MOSFET P-CH 60V 170MA SOT-23"
$ws.Range("H32").Value2 = "DIGIKEY"
$ws.Range("I32").Formula = "=D32*C32"

# Row 33: 47uF ceramic cap
$ws.Range("A33").Value2 = "490-9961-1-ND"
$ws.Range("B33").Value2 = 3
$ws.Range("C33").Value2 = 1
$ws.Range("D33").Value2 = 0.75
$ws.Range("E33").Value2 = "GRM21BR61A476ME15L"
$ws.Range("F33").Value2 = "Murata Electronics North America"
$ws.Range("G33").Value2 = "CAP CER 47UF 10V X5R 0805"
$ws.Range("H33").Value2 = "DIGIKEY"
$ws.Range("I33").Formula = "=D33*C33"

# Row 34: 470pF ceramic cap
$ws.Range("A34").Value2 = "1276-1168-1-ND"
$ws.Range("B34").Value2 = 3
$ws.Range("C34").Value2 = 1
$ws.Range("D34").Value2 = 0.1
$ws.Range("E34").Value2 = "CL10C471JB8NNNC"
$ws.Range("F34").Value2 = "Samsung Electro-Mechanics"
$ws.Range("G34").Value2 = "CAP CER 470PF 50V C0G/NP0 0603"
$ws.Range("H34").Value2 = "DIGIKEY"
$ws.Range("I34").Formula = "=D34*C34"

# Remove the now-duplicated trailing blank row so the sheet keeps the same
# overall extent as before (rows shift down by three, net row count -1)
$ws.Rows("44:44").Delete()

# --- Move the MAX/BOARD total label & formula down onto row 41 ---
$ws.Range("H34").Value2 = ""
$ws.Range("I34").Formula = "=D34*C34"
$ws.Range("H41").Value2 = "MAX/BOARD"
$ws.Range("I41").Formula = "=SUM(I3:I32)"

# --- View/selection tweaks ---
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("B35").Select()
